# Weekly fruit/vegetable price update: the two date-blocks of rows (2-3) and
# (4-5) swap places so the data stays ordered with the newest date on top.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 (was date 44911/Primera) -> becomes date 44915/Especial ---
$ws.Range("D2").Value = 44915
$ws.Range("L2").Value = "Especial"
$ws.Range("M2").Value = 150
$ws.Range("N2").Value = 6000
$ws.Range("O2").Value = 6000
$ws.Range("P2").Value = 6000
$ws.Range("R2").Value = "Provincia de Quillota"
$ws.Range("S2").Value = 1200

# --- Row 3 (was date 44911/Segunda) -> becomes date 44915/Primera ---
$ws.Range("D3").Value = 44915
$ws.Range("L3").Value = "Primera"
$ws.Range("N3").Value = 5000
$ws.Range("O3").Value = 5000
$ws.Range("P3").Value = 5000
$ws.Range("R3").Value = "Provincia de Quillota"
$ws.Range("S3").Value = 1000

# --- Row 4 (was date 44915/Especial) -> becomes date 44911/Primera ---
$ws.Range("D4").Value = 44911
$ws.Range("L4").Value = "Primera"
$ws.Range("M4").Value = 220
$ws.Range("N4").Value = 5000
$ws.Range("O4").Value = 5000
$ws.Range("P4").Value = 5000
$ws.Range("R4").Value = "Región de O'Higgins"
$ws.Range("S4").Value = 1000

# --- Row 5 (was date 44915/Primera) -> becomes date 44911/Segunda ---
$ws.Range("D5").Value = 44911
$ws.Range("L5").Value = "Segunda"
$ws.Range("N5").Value = 4000
$ws.Range("O5").Value = 4000
$ws.Range("P5").Value = 4000
$ws.Range("R5").Value = "Región de O'Higgins"
$ws.Range("S5").Value = 800
